# Sort data rows A2:F43 by Total_profit (column F) descending, keeping header in row 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A1:F43")
$sortRange.Sort($ws.Range("F1"), 2, $null, $null, 1, $null, 1, 1)

# Add three new header cells (G1:I1) matching existing header style
$ws.Range("G1").Value = "Top_3_travel_days"
$ws.Range("H1").Value = "Average_Money_Spent"
$ws.Range("I1").Value = "most_common_languages"
$ws.Range("A1").Copy()
$ws.Range("G1:I1").PasteSpecial(-4122)

# New column data (post-sort row order, rows 2..43)
$gVals = @(
    '[(''Wednesday'', 810), (''Thursday'', 635), (''Sunday'', 631)]',
    '[(''Friday'', 483), (''Thursday'', 482), (''Monday'', 430)]',
    '[(''Monday'', 414), (''Thursday'', 386), (''Friday'', 383)]',
    '[(''Monday'', 552), (''Sunday'', 455), (''Friday'', 363)]',
    '[(''Wednesday'', 440), (''Tuesday'', 371), (''Friday'', 334)]',
    '[(''Monday'', 344), (''Friday'', 318), (''Thursday'', 286)]',
    '[(''Saturday'', 423), (''Friday'', 312), (''Sunday'', 302)]',
    '[(''Thursday'', 377), (''Tuesday'', 335), (''Wednesday'', 313)]',
    '[(''Friday'', 308), (''Saturday'', 255), (''Thursday'', 233)]',
    '[(''Monday'', 236), (''Tuesday'', 202), (''Wednesday'', 173)]',
    '[(''Saturday'', 197), (''Monday'', 193), (''Wednesday'', 184)]',
    '[(''Saturday'', 166), (''Monday'', 157), (''Friday'', 152)]',
    '[(''Monday'', 128), (''Tuesday'', 114), (''Sunday'', 106)]',
    '[(''Saturday'', 195), (''Friday'', 96), (''Tuesday'', 82)]',
    '[(''Friday'', 106), (''Wednesday'', 92), (''Monday'', 71)]',
    '[(''Friday'', 80), (''Sunday'', 72), (''Saturday'', 56)]',
    '[(''Saturday'', 110), (''Friday'', 68), (''Sunday'', 63)]',
    '[(''Friday'', 66), (''Wednesday'', 59), (''Saturday'', 50)]',
    '[(''Friday'', 75), (''Thursday'', 54), (''Monday'', 41)]',
    '[(''Saturday'', 40), (''Sunday'', 40), (''Friday'', 34)]',
    '[(''Sunday'', 57), (''Wednesday'', 43), (''Tuesday'', 37)]',
    '[(''Saturday'', 44), (''Wednesday'', 38), (''Friday'', 36)]',
    '[(''Friday'', 73), (''Saturday'', 60), (''Thursday'', 54)]',
    '[(''Tuesday'', 32), (''Wednesday'', 16), (''Saturday'', 16)]',
    '[(''Friday'', 80), (''Saturday'', 48), (''Sunday'', 21)]',
    '[(''Wednesday'', 40), (''Thursday'', 16), (''Saturday'', 16)]',
    '[(''Saturday'', 47), (''Monday'', 35), (''Wednesday'', 26)]',
    '[(''Wednesday'', 40), (''Thursday'', 40), (''Sunday'', 16)]',
    '[(''Wednesday'', 43), (''Monday'', 31), (''Friday'', 18)]',
    '[(''Friday'', 56), (''Saturday'', 35), (''Tuesday'', 28)]',
    '[(''Thursday'', 26), (''Sunday'', 25), (''Saturday'', 13)]',
    '[(''Saturday'', 32), (''Friday'', 32), (''Tuesday'', 16)]',
    '[(''Sunday'', 22), (''Monday'', 22), (''Thursday'', 17)]',
    '[(''Tuesday'', 24), (''Friday'', 24), (''Thursday'', 16)]',
    '[(''Thursday'', 20), (''Tuesday'', 15), (''Monday'', 10)]',
    '[(''Thursday'', 25), (''Wednesday'', 11), (''Friday'', 8)]',
    '[(''Tuesday'', 24), (''Friday'', 16), (''Wednesday'', 8)]',
    '[(''Saturday'', 11), (''Thursday'', 10), (''Friday'', 8)]',
    '[(''Tuesday'', 14), (''Sunday'', 10), (''Wednesday'', 8)]',
    '[(''Sunday'', 16), (''Wednesday'', 8), (''Thursday'', 8)]',
    '[(''Saturday'', 22), (''Thursday'', 10), (''Wednesday'', 5)]',
    '[(''Sunday'', 10), (''Wednesday'', 10), (''Monday'', 10)]'
)

$hVals = @(
    104.7646530902202,
    97.91842900302115,
    98.8620759289176,
    95.22041946972695,
    99.91595216191352,
    98.68495762711865,
    96.90059494298463,
    87.53344173441735,
    99.51302403204274,
    93.98205223880596,
    89.30600000000001,
    84.68898043254377,
    91.57888730385164,
    79.6590909090909,
    95.16538461538461,
    121.9487179487179,
    103.6697247706422,
    84.58237082066869,
    97.1888030888031,
    111.1859821428571,
    94.75728155339806,
    97.52226804123711,
    61.03460992907801,
    174,
    75.80714285714286,
    152.9090909090909,
    79.57888324873097,
    116,
    86.19083969465649,
    62.77844311377245,
    109.1075862068966,
    100.3846153846154,
    97.99654320987655,
    109.5555555555556,
    95.15753086419754,
    115.4038461538462,
    95.8925,
    87.73333333333333,
    90.2,
    87,
    59.01666666666667,
    80.88888888888889
)

$iVals = @(
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'', ''German'']',
    '[''English'']',
    '[''English'']',
    '[''French'']',
    '[''English'']',
    '[''German'']',
    '[''English'']',
    '[''Spanish'']',
    '[''German'']',
    '[''Spanish'']',
    '[''English'', ''German'']',
    '[''English'']',
    '[''Italian'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''English'']',
    '[''Spanish'']',
    '[''Spanish'']',
    '[''English'']',
    '[''German'']',
    '[''Spanish'']',
    '[''English'']'
)

for ($idx = 0; $idx -lt 42; $idx++) {
    $r = $idx + 2
    $ws.Cells.Item($r, 7).Value = $gVals[$idx]
    $ws.Cells.Item($r, 8).Value = $hVals[$idx]
    $ws.Cells.Item($r, 9).Value = $iVals[$idx]
}

Write-Host "Edit complete"
